$d = $word.ActiveDocument

# The last paragraph of the document currently reads:
#   "If you pass in a value that is not a string it will call to_s ||non_string||"
# and carries the (hidden) _GoBack bookmark right at its very end, just before
# the paragraph mark. We need to:
#   1. Add a whole new paragraph right after it.
#   2. Fill that new paragraph with the new text about un-matched placeholders.
#   3. Move the _GoBack bookmark so that it sits in the middle of that new
#      paragraph's text (between "...||stay_on_the_page" and "|| so that...").

# Step 1: Remove the existing _GoBack bookmark. We will recreate it later, in
# its new location, once the surrounding text exists.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# Step 2: Locate the "non_string" paragraph (the last real paragraph of the
# body, right before the trailing empty paragraph / sectPr) and split a new,
# empty paragraph in right after it.
# NOTE: Paragraph.Index is unreliable in this runtime, so we track the
# 1-based position with our own loop counter instead of trusting .Index.
$nonStringParaNumber = -1
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*non_string*") {
        $nonStringParaNumber = $i
        break
    }
}

$nonStringParagraph = $d.Paragraphs.Item($nonStringParaNumber)
$endOfParagraph = $nonStringParagraph.Range
$endOfParagraph.Collapse(0)   # wdCollapseEnd
$endOfParagraph.InsertParagraphAfter()

# Step 3: Insert all of the new paragraph's text in one shot, right after the
# paragraph mark we just created.
$newParagraph = $d.Paragraphs.Item($nonStringParaNumber + 1)
$insertionPoint = $d.Range($newParagraph.Range.Start, $newParagraph.Range.Start)

$firstPart = "If there is a placeholder that doesn" + [char]0x2019 + "t match a data provider"
$secondPart = " key it should ||stay_on_the_page"
$thirdPart = "|| so that you can debug. Even if a value is nil, it should have the key in the data provider."

$insertionPoint.InsertAfter($firstPart + $secondPart + $thirdPart)

# Step 4: Re-add the _GoBack bookmark as a zero-length bookmark sitting right
# between "||stay_on_the_page" and "|| so that...".
$bookmarkPos = $newParagraph.Range.Start + $firstPart.Length + $secondPart.Length
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
